# "Repayment schedule" becomes the active/selected sheet (was "Edit Repayment
# Schedule"), a new blank column is inserted before column N ("Late"), and
# the selection on that sheet moves to J14.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column at N, shifting the existing N/O/P ("Late",
# "heading", "Outstanding") columns one to the right.
$ws.Columns("N:N").Insert() | Out-Null

# Excel's ColumnWidth property is expressed in character widths and gets an
# implicit ~0.8333 padding when stored as the raw OOXML "width" attribute, so
# back it off to land on the target raw width of 11.
$ws.Columns("N:N").ColumnWidth = 10.166666666666666

# Move the selection to J14, matching the saved view state.
$ws.Range("J14").Select() | Out-Null
